$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Swap column E (codeforiati:category-name) and column F (codeforiati:group-code)
# for every row, including the header. Use .Copy() (instead of .Value = ...) so
# text-typed cells (e.g. numeric-looking codes like "110") keep their original
# string type instead of being re-interpreted as numbers, and so no new cell
# style gets introduced.
$scratch = $ws.Range("Z1")
for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eCell.Copy($scratch)
    $fCell.Copy($eCell)
    $scratch.Copy($fCell)
}
$scratch.ClearContents()
